$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates. These cells are text-typed (string representations of
# numbers) and must keep their exact textual formatting (e.g. trailing zeros), so
# we force the cells to Text number format before writing the new value - this
# prevents Excel from re-interpreting the string as a floating point number.
$priceUpdates = @{
    2  = "281.80"
    3  = "20.62"
    4  = "6.236"
    5  = "0.06149"
    7  = "6.568"
    8  = "1.499"
    9  = "0.8177"
    10 = "0.01384"
    11 = "0.1635"
    12 = "0.08357"
    13 = "0.03544"
    14 = "0.03188"
    15 = "0.09137"
    16 = "3.706"
    17 = "0.001644"
    18 = "0.04707"
    19 = "0.006444"
    20 = "0.006160"
    21 = "0.001069"
    23 = "3.769"
    25 = "0.3356"
    40 = "0.04682"
    41 = "0.007176"
    42 = "0.1101"
    43 = "0.003483"
    45 = "0.00006635"
    48 = "0.002955"
    49 = "0.00001902"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# Rows 42 and 43 swap coin identity (CEJI <-> BKEXToken), including their
# Coin name, Link and Volume(1h) label columns.
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"

Write-Host "Applied price/symbol list updates"
